# Apply cryptocurrency price/volume updates to Sheet1 (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the "Price" cells that would otherwise be auto-recognized as
# numbers (single decimal point) so the assigned value is stored as text,
# matching the original inline-string cell type for column D.
$textPriceCells = "D5","D7","D10","D11","D14","D17","D20","D22","D25","D28","D30","D37","D38","D42","D43","D45","D47"
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.319.31"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "1.620.71"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "212.11"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.484"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "18.79"
$ws.Range("E10").Value = "  +4.58%  "
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "1.846.51"
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("D13").Value = "1.621.74"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").Value = "26.313.33"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "62.30"
$ws.Range("E17").Value = "  +3.57%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "201.72"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").Value = "9.34"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("E24").Value = "  -3.55%  "
$ws.Range("D25").Value = "144.85"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").Value = "15.18"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("D30").Value = "0.0517"
$ws.Range("E30").Value = "  +9.03%  "
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("E35").Value = "  +2.55%  "
$ws.Range("D36").Value = "1.179.90"
$ws.Range("E36").Value = "  +4.79%  "
$ws.Range("D37").Value = "0.0164"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "0.809"
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  +1.31%  "
$ws.Range("D42").Value = "0.786"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("D43").Value = "5.35"
$ws.Range("E43").Value = "  +4.46%  "
$ws.Range("D44").Value = "1.757.66"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("D45").Value = "92.73"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("E46").Value = "  +3.25%  "
$ws.Range("D47").Value = "53.78"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("E51").Value = "  +2.14%  "
